$d = $word.ActiveDocument

# --- Add the three new character styles (order matters for OOXML output) ---

$GaNStyle = $d.Styles.Add("GaNStyle", 2)
$GaNStyle.Font.Name = "Calibri"
$GaNStyle.Font.Size = 14

$GaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$GaNParagraph.Font.Name = "Calibri"
$GaNParagraph.Font.Size = 10

$GaNLinks = $d.Styles.Add("GaNLinks", 2)
$GaNLinks.Font.Name = "Calibri"
$GaNLinks.Font.Size = 9.5
$GaNLinks.Font.Bold = $true
$GaNLinks.Font.Color = 8388608
$GaNLinks.Font.Underline = 1

# --- Apply GaNStyle to every "2022 Fechas de la campaña..." run (4 occurrences) ---

$fechasText = "2022 Fechas de la campaña para Constelación de botas: 14-23 de mayo, 13-22 de junio, 12-21 de julio"
$rng = $d.Content
while ($rng.Find.Execute($fechasText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNStyle"
}

# --- Apply GaNParagraph to the intro paragraph run ---

$parrafoText = "Usted está participando en una campaña mundial para observar y registrar las estrellas visibles más débiles como un medio para medir la contaminación lumínica en un lugar determinado. Localizando y observando la  Constelación de botas en el cielo nocturno y comparándolo con las cartas estelares, la gente de todo el mundo aprenderán cómo las luces de su comunidad contribuyen a la contaminación lumínica. Sus contribuciones a la base de datos en línea documentarán el cielo nocturno visible."
$rng = $d.Content
if ($rng.Find.Execute($parrafoText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNParagraph"
}

# --- Apply GaNLinks to the map-link run ---

$linkText = "(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
$rng = $d.Content
if ($rng.Find.Execute($linkText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNLinks"
}
